$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows for 2000年-2009年 (rows 2-11). This shifts the
# 2010年-2020年 data (previously rows 12-22) up to rows 2-12.
$ws.Range("A2:F11").Delete(-4162) | Out-Null

# Append the new data rows for 2021年 and 2022年, copying the
# formatting used by the other "year" cells in column A.
$ws.Range("A13").Value = "2021年"
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A14").Value = "2022年"
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("B13").Value = 100
$ws.Range("C13").Value = 0.342430221903317
$ws.Range("D13").Value = 23.2714861268005
$ws.Range("E13").Value = 0.77
$ws.Range("F13").Value = 7.72

$ws.Range("B14").Value = 100
$ws.Range("E14").Value = 0.67
